# Fuel Prod Imp Exp Balancing Priorities.xlsx - "updated 4.0 files and mdl"
#
# Semantic edits captured in the canonical-XML diff:
#   1. About!C1   - the "as-of" date stamp moves from 1/3/2024 (45294)
#                    to 3/28/2024 (45379).
#   2. FPIEBP!B3:D3 - "hard coal" production/imports/exports priority
#                    values are reshuffled from (3,2,1) to (1,3,2).
#   3. FPIEBP sheet's selected/active cell moves from F4 to E3 (the
#                    sheet stays the active/selected tab throughout).
#
# (Everything else in the raw diff - fileVersion/rupBuild, the
# revisionPtr GUID, workbookView window geometry, the theme's cosmetic
# "Office" -> "Office 2013 - 2022" rename, and the sub-point-wide
# column/row metric jitter - is Excel build/session bookkeeping that
# isn't data and has no corresponding Excel object model call, so it's
# intentionally left alone here.)

$wb = $excel.ActiveWorkbook

# --- About sheet: update the date stamp in C1 ---------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- FPIEBP sheet: reorder the "hard coal" priority row -----------------
$wsFpiebp = $wb.Worksheets.Item("FPIEBP")
$wsFpiebp.Range("B3").Value = 1
$wsFpiebp.Range("C3").Value = 3
$wsFpiebp.Range("D3").Value = 2

# Move the active selection to E3 (FPIEBP is already the active sheet,
# so this doesn't disturb ActiveTab / tabSelected elsewhere).
$wsFpiebp.Range("E3").Select()
